$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8: requirement 4 "Concealed pieces can capture other pieces and put
#     King in check." flips from In-Progress (yellow) to Done (green). ---
$ws.Range("A8").Interior.Color = 65280
$ws.Range("B8").Interior.Color = 65280
$ws.Range("C8").Interior.Color = 65280
$ws.Range("C8").Value = "Done"
$ws.Range("D8").Interior.Color = 65280

# --- Row 10: requirement 4.2 "Concealed pieces can put the King in check."
#     flips from In-Progress (yellow) to Done (green). ---
$ws.Range("A10").Interior.Color = 65280
$ws.Range("B10").Interior.Color = 65280
$ws.Range("C10").Interior.Color = 65280
$ws.Range("C10").Value = "Done"
$ws.Range("D10").Interior.Color = 65280

# --- Row 18: requirement 6 "Castling move has been made unavailable due to
#     rule changes." flips from In-Progress (yellow) to Done (green). ---
$ws.Range("A18").Interior.Color = 65280
$ws.Range("B18").Interior.Color = 65280
$ws.Range("C18").Interior.Color = 65280
$ws.Range("C18").Value = "Done"
$ws.Range("D18").Interior.Color = 65280

# --- Update the saved selection/active cell shown when the sheet re-opens. ---
$ws.Range("G21:G22").Select()
